$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3:22 down to 4:23
$ws.Rows(3).Insert()

# Populate the new row 3 with its data (same constant columns as the rest
# of the dataset, plus the new record's own values)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 45043
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101007
$ws.Range("J3").Value = "Kiwi"
$ws.Range("K3").Value = "Hayward"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21500
$ws.Range("Q3").Value = "`$/bandeja 18 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1194
$ws.Range("T3").Value = 18

# D column keeps the date number format used by the other rows
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
